# Update cryptos worksheet with latest crypto price/volume data
# (values matching diff: mirrors an automated "Updated cryptos list" commit)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.849.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.79%  "
$ws.Range("D3").Value = "'1.756.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "'327.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.80%  "
$ws.Range("E6").Value = "  +0.32%  "
$ws.Range("D7").Value = "'0.4687"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.19%  "
$ws.Range("D8").Value = "'0.3499"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.20%  "
$ws.Range("D9").Value = "'41.97"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.74%  "
$ws.Range("D10").Value = "'0.07357"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.53%  "
$ws.Range("D11").Value = "'1.081"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.74%  "
$ws.Range("D12").Value = "'1.000"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("E13").Value = "  -1.78%  "
$ws.Range("D14").Value = "'5.980"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.75%  "
$ws.Range("D15").Value = "'7.146"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.82%  "
$ws.Range("D16").Value = "'1.755.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.70%  "
$ws.Range("D17").Value = "'91.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.96%  "
$ws.Range("D18").Value = "'0.00001053"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("D19").Value = "'0.06402"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.61%  "
$ws.Range("E20").Value = "  +0.31%  "
$ws.Range("D21").Value = "'16.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.09%  "
$ws.Range("D22").Value = "'5.749"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.43%  "
$ws.Range("D23").Value = "'27.876.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("D24").Value = "'11.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.05%  "
$ws.Range("D25").Value = "'2.147"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.56%  "
$ws.Range("D26").Value = "'161.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.35%  "
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("D28").Value = "'1.958.56"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.46%  "
$ws.Range("D29").Value = "'2.148"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.81%  "
$ws.Range("D30").Value = "'122.55"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.86%  "
$ws.Range("D31").Value = "'1.067"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.21%  "
$ws.Range("D32").Value = "'0.09331"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.46%  "
$ws.Range("E33").Value = "  -0.50%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  -1.04%  "
$ws.Range("D36").Value = "'11.61"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.83%  "
$ws.Range("D37").Value = "'0.06062"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").Value = "'0.2063"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.38%  "
$ws.Range("D39").Value = "'4.889"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.19%  "
$ws.Range("D40").Value = "'0.6128"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.71%  "
$ws.Range("D41").Value = "'1.178"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.62%  "
$ws.Range("D42").Value = "'7.761"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.39%  "
$ws.Range("D43").Value = "'1.349"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.68%  "
$ws.Range("D44").Value = "'13.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.40%  "
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("D46").Value = "'0.5768"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.93%  "
$ws.Range("D47").Value = "'122.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("E48").Value = "  -1.31%  "
$ws.Range("D49").Value = "'0.06794"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.08%  "
$ws.Range("E50").Value = "  -1.40%  "
$ws.Range("D51").Value = "'72.08"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.69%  "
